# Add season-record columns (Wins, Losses, Ties) to the roster table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: copy the style from the last existing header cell (AC1, style index 1:
# bold, centered, thin border) onto the three new header cells, then set their text.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Data rows 2-42: every player row gets the team's season record.
for ($r = 2; $r -le 42; $r++) {
    $ws.Cells.Item($r, 30).Value = 86  # AD = Wins
    $ws.Cells.Item($r, 31).Value = 76  # AE = Losses
    $ws.Cells.Item($r, 32).Value = 0   # AF = Ties
}
